$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add new row 3 data to the Logs sheet
$logs.Range("A3").Value = "Kun jij dit even regelen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D3").Value = "Overig"
$logs.Range("F3").Value = "2025-07-29 21:28:46"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# Update dashboard count
$dash.Range("B2").Value = 2

# Extend conditional formatting ranges to include row 3
foreach ($col in @("D","G","H","I","J")) {
    $oldRange = $logs.Range($col + "2")
    $newRange = $logs.Range($col + "2:" + $col + "3")
    foreach ($fc in $oldRange.FormatConditions) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

